# Applies the LOUISIANA_2017 cleaning edit:
#  1. Rename header columns to snake_case English names
#  2. Title-case the Spanish articles/prepositions "de", "del", "las", "los", "el"
#     (but NOT "la") when they appear as a non-leading word within a municipality /
#     state name in columns A and B
#  3. Refresh the percentage column D values (tiny float re-computation: 3 / 3176)
#  4. Remove the trailing footnote rows (812-816), shrinking the used range to A1:D810

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row renames
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,1).Value2 = "mx_state"
$ws.Cells.Item(1,2).Value2 = "mx_municipality"
$ws.Cells.Item(1,3).Value2 = "n_matriculas"
$ws.Cells.Item(1,4).Value2 = "pct_matriculas"

# ---------------------------------------------------------------------------
# 2. Title-case the articles/prepositions in columns A and B (rows 2-810)
# ---------------------------------------------------------------------------
$articles = @("de", "del", "las", "los", "el")

for ($r = 2; $r -le 810; $r++) {
    foreach ($col in 1,2) {
        $cell = $ws.Cells.Item($r, $col)
        $text = $cell.Value2
        if ($text -ne $null -and ($text -is [string])) {
            $words = $text -split ' '
            $changed = $false
            for ($i = 1; $i -lt $words.Length; $i++) {
                $w = $words[$i]
                if ($articles -contains $w) {
                    $words[$i] = $w.Substring(0,1).ToUpper() + $w.Substring(1)
                    $changed = $true
                }
            }
            if ($changed) {
                $cell.Value2 = ($words -join ' ')
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 3. Refresh column D float values (3 / 3176 got re-computed to the adjacent
#    float64 value during the data regeneration)
# ---------------------------------------------------------------------------
$oldVal = 0.0009445843828715365
$newVal = 0.0009445843828715364

for ($r = 2; $r -le 810; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $v = $cell.Value2
    if ($v -ne $null -and -not ($v -is [string])) {
        if ($v -eq $oldVal) {
            $cell.Value2 = $newVal
        }
    }
}

# ---------------------------------------------------------------------------
# 4. Drop the trailing footnote rows (812-816); row 811 is already blank
# ---------------------------------------------------------------------------
$ws.Range("A812:A816").EntireRow.Delete() | Out-Null
